# Fruta / hortaliza, semanal
# Insert a new weekly record as row 428, pushing the existing rows
# (428-487) down to (429-488).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 428, shifting everything
# below it (including the dimension) down by one row.
$ws.Rows(428).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(428, 1).Value = 8
$ws.Cells.Item(428, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(428, 3).Value = "Coquimbo"
$ws.Cells.Item(428, 4).Value = 45124
$ws.Cells.Item(428, 5).Value = 4
$ws.Cells.Item(428, 6).Value = 100112021
$ws.Cells.Item(428, 7).Value = "Ají"
$ws.Cells.Item(428, 8).Value = "Inferno"
$ws.Cells.Item(428, 9).Value = "Primera"
$ws.Cells.Item(428, 10).Value = 360
$ws.Cells.Item(428, 11).Value = 13500
$ws.Cells.Item(428, 12).Value = 14000
$ws.Cells.Item(428, 13).Value = 13750
$ws.Cells.Item(428, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(428, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(428, 16).Value = 1375
$ws.Cells.Item(428, 17).Value = 10
$ws.Cells.Item(428, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(428, 4).NumberFormat = $ws.Cells.Item(429, 4).NumberFormat
